# Adiciona CS_Consorcio_Contemplado.py ao processo de extracao
# -> atualiza os valores extraidos (carta, entrada, parcelas, consorcio,
#    fluxo de pagamento) e acrescenta o novo registro VC1020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (VC1001)
$ws.Range("C2").Value = "39100,00"
$ws.Range("D2").Value = "24855,00"
$ws.Range("E2").Value = "33"
$ws.Range("H2").Value = "33 x R$ 538.97"

# Row 3 (VC1002)
$ws.Range("C3").Value = "39300,00"
$ws.Range("D3").Value = "29965,00"
$ws.Range("E3").Value = "44"
$ws.Range("F3").Value = "Bradesco"
$ws.Range("H3").Value = "44 x R$ 344.95"

# Row 4 (VC1003)
$ws.Range("C4").Value = "43200,00"
$ws.Range("D4").Value = "31460,00"
$ws.Range("E4").Value = "39"
$ws.Range("H4").Value = "39 x R$ 318.46"

# Row 5 (VC1004)
$ws.Range("C5").Value = "44700,00"
$ws.Range("D5").Value = "25735,00"
$ws.Range("E5").Value = "26"
$ws.Range("F5").Value = "Porto Seguro"
$ws.Range("H5").Value = "26 x R$ 1057.23"

# Row 6 (VC1005)
$ws.Range("C6").Value = "51000,00"
$ws.Range("D6").Value = "29450,00"
$ws.Range("E6").Value = "37"
$ws.Range("H6").Value = "37 x R$ 1044.02"

# Row 7 (VC1006)
$ws.Range("C7").Value = "59350,00"
$ws.Range("D7").Value = "33867,50"
$ws.Range("E7").Value = "45"
$ws.Range("F7").Value = "Bradesco"
$ws.Range("H7").Value = "45 x R$ 718.11"

# Row 8 (VC1007)
$ws.Range("C8").Value = "62400,00"
$ws.Range("D8").Value = "48020,00"
$ws.Range("E8").Value = "32"
$ws.Range("F8").Value = "Porto Seguro"
$ws.Range("H8").Value = "32 x R$ 546.01"

# Row 9 (VC1008)
$ws.Range("C9").Value = "62600,00"
$ws.Range("D9").Value = "25130,00"
$ws.Range("E9").Value = "45"
$ws.Range("H9").Value = "45 x R$ 1549.77"

# Row 10 (VC1009)
$ws.Range("C10").Value = "64100,00"
$ws.Range("D10").Value = "40105,00"
$ws.Range("E10").Value = "34"
$ws.Range("F10").Value = "Porto Seguro"
$ws.Range("H10").Value = "34 x R$ 1228.00"

# Row 11 (VC1010)
$ws.Range("C11").Value = "65800,00"
$ws.Range("D11").Value = "32190,00"
$ws.Range("E11").Value = "60"
$ws.Range("F11").Value = "Porto Seguro"
$ws.Range("H11").Value = "60 x R$ 1145.00"

# Row 12 (VC1011)
$ws.Range("C12").Value = "71700,00"
$ws.Range("D12").Value = "36485,00"
$ws.Range("E12").Value = "17"
$ws.Range("F12").Value = "Porto Seguro"
$ws.Range("H12").Value = "17 x R$ 2480.34"

# Row 13 (VC1012)
$ws.Range("C13").Value = "77600,00"
$ws.Range("D13").Value = "47780,00"
$ws.Range("E13").Value = "39"
$ws.Range("H13").Value = "39 x R$ 1022.32"

# Row 14 (VC1013)
$ws.Range("C14").Value = "89000,00"
$ws.Range("D14").Value = "46450,00"
$ws.Range("E14").Value = "20"
$ws.Range("H14").Value = "20 x R$ 3184.65"

# Row 15 (VC1014)
$ws.Range("C15").Value = "92150,00"
$ws.Range("D15").Value = "52107,50"
$ws.Range("E15").Value = "67"
$ws.Range("H15").Value = "67 x R$ 1015.55"

# Row 16 (VC1015)
$ws.Range("C16").Value = "116500,00"
$ws.Range("D16").Value = "57825,00"
$ws.Range("E16").Value = "35"
$ws.Range("H16").Value = "35 x R$ 2729.07"

# Row 17 (VC1016)
$ws.Range("C17").Value = "116500,00"
$ws.Range("D17").Value = "55725,00"
$ws.Range("E17").Value = "30"
$ws.Range("H17").Value = "30 x R$ 3193.82"

# Row 18 (VC1017) - muda de Veiculos para Imoveis
$ws.Range("B18").Value = "Imóveis"
$ws.Range("C18").Value = "99600,00"
$ws.Range("D18").Value = "54880,00"
$ws.Range("E18").Value = "100"
$ws.Range("H18").Value = "100 x R$ 616.00"

# Row 19 (VC1018)
$ws.Range("C19").Value = "140000,00"
$ws.Range("D19").Value = "84000,00"
$ws.Range("E19").Value = "60"
$ws.Range("H19").Value = "60 x R$ 1356.82"

# Row 20 (VC1019)
$ws.Range("C20").Value = "404000,00"
$ws.Range("D20").Value = "247200,00"
$ws.Range("E20").Value = "158"
$ws.Range("H20").Value = "158 x R$ 1898.72"

# Nova linha 21 (VC1020) - copia a formatacao da linha 20 e preenche os dados
$ws.Range("A20:J20").Copy()
$ws.Range("A21:J21").PasteSpecial(-4122)
$ws.Range("A21").Value = "VC1020"
$ws.Range("B21").Value = "Imóveis"
$ws.Range("C21").Value = "544000,00"
$ws.Range("D21").Value = "329200,00"
$ws.Range("E21").Value = "158"
$ws.Range("F21").Value = "Porto Seguro"
$ws.Range("G21").Value = "Disponível"
$ws.Range("H21").Value = "60 x R$ 3255.00`n98 x R$ 1898.00"
$ws.Rows(21).AutoFit()

# Coluna H (Fluxo de Pagamento) fica um pouco mais larga (32 -> 33)
$ws.Columns("H").ColumnWidth = 32.15
